$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = @(66.47695399999999, 199.430862, 0.04311983106164722, 0.04311983106164721, 3.0, 1.0, 0.3708513333333334, 1.112554, 0.07446453063452614, 0.07446453063452615, 24.65306702683867, 221.877603241548, 0.003210897981045621, 0.003210897981045621)
$row3 = @(66.47695399999999, 199.430862, 0.04311983106164722, 0.04311983106164721, 3.0, 1.0, 2.230730666666667, 6.692192, 0.4479161786269527, 0.4479161786269528, 148.2921799143893, 1334.629619229504, 0.0193140699521728, 0.0193140699521728)
$row4 = @(66.47695399999999, 199.430862, 0.04311983106164722, 0.04311983106164721, 3.0, 1.0, 1.029162, 3.087486, 0.2066490218278579, 0.2066490218278579, 68.415554932548, 615.739994392932, 0.008910670910271881, 0.00891067091027188)
$row5 = @(66.47695399999999, 199.430862, 0.04311983106164722, 0.04311983106164721, 3.0, 1.0, 0.8350726666666667, 2.505218, 0.1676771487111334, 0.1676771487111334, 55.51308724865733, 499.617785237916, 0.00723021032532277, 0.00723021032532277)
$row6 = @(66.47695399999999, 199.430862, 0.04311983106164722, 0.04311983106164721, 3.0, 1.0, 0.5144246666666666, 1.543274, 0.1032931201995298, 0.1032931201995298, 34.19738490246533, 307.776464122188, 0.004453981892834145, 0.004453981892834145)
$row7 = @(1361.379069, 4084.137207, 0.8830494168872806, 0.8830494168872804, 3.0, 1.0, 0.3708513333333334, 1.112554, 0.07446453063452614, 0.07446453063452615, 504.869242910742, 4543.823186196678, 0.06575586035560335, 0.06575586035560335)
$row8 = @(1361.379069, 4084.137207, 0.8830494168872806, 0.8830494168872804, 3.0, 1.0, 2.230730666666667, 6.692192, 0.4479161786269527, 0.4479161786269528, 3036.870038176416, 27331.83034358774, 0.3955321203509096, 0.3955321203509096)
$row9 = @(1361.379069, 4084.137207, 0.8830494168872806, 0.8830494168872804, 3.0, 1.0, 1.029162, 3.087486, 0.2066490218278579, 0.2066490218278579, 1401.079605410178, 12609.7164486916, 0.1824812982254168, 0.1824812982254168)
$row10 = @(1361.379069, 4084.137207, 0.8830494168872806, 0.8830494168872804, 3.0, 1.0, 0.8350726666666667, 2.505218, 0.1676771487111334, 0.1676771487111334, 1136.850449494014, 10231.65404544613, 0.1480672083946882, 0.1480672083946882)
$row11 = @(1361.379069, 4084.137207, 0.8830494168872806, 0.8830494168872804, 3.0, 1.0, 0.5144246666666666, 1.543274, 0.1032931201995298, 0.1032931201995298, 700.3269737773019, 6302.942763995718, 0.09121292956066258, 0.09121292956066258)
$row12 = @(44.831112, 134.493336, 0.02907940059566787, 0.02907940059566786, 3.0, 1.0, 0.3708513333333334, 1.112554, 0.07446453063452614, 0.07446453063452615, 16.625677660016, 149.631098940144, 0.002165383916489768, 0.002165383916489768)
$row13 = @(44.831112, 134.493336, 0.02907940059566787, 0.02907940059566786, 3.0, 1.0, 2.230730666666667, 6.692192, 0.4479161786269527, 0.4479161786269528, 100.006136359168, 900.055227232512, 0.01302513399157389, 0.01302513399157388)
$row14 = @(44.831112, 134.493336, 0.02907940059566787, 0.02907940059566786, 3.0, 1.0, 1.029162, 3.087486, 0.2066490218278579, 0.2066490218278579, 46.13847688814401, 415.246291993296, 0.006009229688435193, 0.006009229688435192)
$row15 = @(44.831112, 134.493336, 0.02907940059566787, 0.02907940059566786, 3.0, 1.0, 0.8350726666666667, 2.505218, 0.1676771487111334, 0.1676771487111334, 37.437236247472, 336.935126227248, 0.004875950978110423, 0.004875950978110422)
$row16 = @(44.831112, 134.493336, 0.02907940059566787, 0.02907940059566786, 3.0, 1.0, 0.5144246666666666, 1.543274, 0.1032931201995298, 0.1032931201995298, 23.062229846896, 207.560068622064, 0.0030037020210586, 0.0030037020210586)
$row17 = @(52.83062100000001, 158.491863, 0.0342682285413064, 0.03426822854130639, 3.0, 1.0, 0.3708513333333334, 1.112554, 0.07446453063452614, 0.07446453063452615, 19.59230623867801, 176.330756148102, 0.002551767554005054, 0.002551767554005054)
$row18 = @(52.83062100000001, 158.491863, 0.0342682285413064, 0.03426822854130639, 3.0, 1.0, 2.230730666666667, 6.692192, 0.4479161786269527, 0.4479161786269528, 117.850886403744, 1060.657977633696, 0.01534929397653704, 0.01534929397653704)
$row19 = @(52.83062100000001, 158.491863, 0.0342682285413064, 0.03426822854130639, 3.0, 1.0, 1.029162, 3.087486, 0.2066490218278579, 0.2066490218278579, 54.37126756960202, 489.3414081264181, 0.007081495907834448, 0.007081495907834448)
$row20 = @(52.83062100000001, 158.491863, 0.0342682285413064, 0.03426822854130639, 3.0, 1.0, 0.8350726666666667, 2.505218, 0.1676771487111334, 0.1676771487111334, 44.117407560126, 397.0566680411341, 0.00574599885318774, 0.005745998853187739)
$row21 = @(52.83062100000001, 158.491863, 0.0342682285413064, 0.03426822854130639, 3.0, 1.0, 0.5144246666666666, 1.543274, 0.1032931201995298, 0.1032931201995298, 27.177374597718, 244.596371379462, 0.00353967224974212, 0.00353967224974212)
$row22 = @(16.16161433333333, 48.484843, 0.01048312291409786, 0.01048312291409786, 3.0, 1.0, 0.3708513333333334, 1.112554, 0.07446453063452614, 0.07446453063452615, 5.993556224335778, 53.942006019022, 0.0007806208273823434, 0.0007806208273823434)
$row23 = @(16.16161433333333, 48.484843, 0.01048312291409786, 0.01048312291409786, 3.0, 1.0, 2.230730666666667, 6.692192, 0.4479161786269527, 0.4479161786269528, 36.05220871620622, 324.469878445856, 0.00469556035575936, 0.00469556035575936)
$row24 = @(16.16161433333333, 48.484843, 0.01048312291409786, 0.01048312291409786, 3.0, 1.0, 1.029162, 3.087486, 0.2066490218278579, 0.2066490218278579, 16.632919330522, 149.696273974698, 0.002166327095899527, 0.002166327095899526)
$row25 = @(16.16161433333333, 48.484843, 0.01048312291409786, 0.01048312291409786, 3.0, 1.0, 0.8350726666666667, 2.505218, 0.1676771487111334, 0.1676771487111334, 13.49612237897489, 121.465101410774, 0.001757780159824278, 0.001757780159824278)
$row26 = @(16.16161433333333, 48.484843, 0.01048312291409786, 0.01048312291409786, 3.0, 1.0, 0.5144246666666666, 1.543274, 0.1032931201995298, 0.1032931201995298, 8.31393306622022, 74.825397595982, 0.001082834475232356, 0.001082834475232356)

$allRows = @($row2, $row3, $row4, $row5, $row6, $row7, $row8, $row9, $row10, $row11, $row12, $row13, $row14, $row15, $row16, $row17, $row18, $row19, $row20, $row21, $row22, $row23, $row24, $row25, $row26)
$rowNums = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26)

for ($i = 0; $i -lt $rowNums.Count; $i++) {
    $r = $rowNums[$i]
    $vals = $allRows[$i]
    for ($j = 0; $j -lt $vals.Count; $j++) {
        $c = 7 + $j
        $ws.Cells.Item($r, $c).Value = $vals[$j]
    }
}

Write-Output "Update complete"